$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.581.97"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.531.90"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.96"
$ws.Range("E5").Value = "  +4.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.87"
$ws.Range("E6").Value = "  -3.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  -2.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.31"
$ws.Range("E10").Value = "  -1.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.70"
$ws.Range("E12").Value = "  -0.84%  "

$ws.Range("E13").Value = "  -0.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.920.30"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.89"
$ws.Range("E15").Value = "  +4.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.533.96"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.865"
$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.644.60"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.06"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("E20").Value = "  +0.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0967"
$ws.Range("E21").Value = "  -2.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.97"
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.80"
$ws.Range("E23").Value = "  -0.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("E25").Value = "  -3.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.98"
$ws.Range("E26").Value = "  -3.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.37"
$ws.Range("E28").Value = "  +2.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.58"
$ws.Range("E29").Value = "  +2.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.22"
$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.97"
$ws.Range("E31").Value = "  -3.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.85"
$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.14"
$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.32"
$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.90"
$ws.Range("E35").Value = "  -1.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0787"
$ws.Range("E36").Value = "  -1.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.62"
$ws.Range("E37").Value = "  -0.73%  "

$ws.Range("E38").Value = "  -3.04%  "

$ws.Range("E39").Value = "  -0.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.86"
$ws.Range("E40").Value = "  -4.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.31"
$ws.Range("E41").Value = "  +9.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.84"
$ws.Range("E42").Value = "  -1.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.35"
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0300"
$ws.Range("E45").Value = "  -1.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.028.73"
$ws.Range("E46").Value = "  -3.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.53"
$ws.Range("E47").Value = "  -2.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.89"
$ws.Range("E48").Value = "  -1.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.775.93"
$ws.Range("E49").Value = "  -0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.03"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.190"
$ws.Range("E51").Value = "  -1.26%  "
